# Applies the "calibration page & time set & elapsed time" edit:
#  - Typography sheet: register a new font entry "Typography_05" in row 12
#    (verdana.ttf, size 14, bpp 4, fallback "?"), mirroring rows 4-11.
#  - Translation sheet: renumber the text-id placeholders for rows 68-70,
#    and populate rows 71-105 with the calibration-page / time-set /
#    elapsed-time / numeric-keypad-value text entries.

$wb = $excel.ActiveWorkbook

# Helper: force a value to be written as literal TEXT (shared string),
# even when it looks like a number (e.g. "0000"), and without acquiring
# a custom number-format style along the way. We do this by writing a
# formula that evaluates to the literal string, then converting the
# whole touched range to static values via Copy + PasteSpecial(Values).
function Set-TextCell($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
}
# --- Sheet "Typography": add new font row 12 (Typography_05) ---
$ws1 = $wb.Worksheets.Item("Typography")

Set-TextCell $ws1.Cells.Item(12, 2) 'Typography_05'
Set-TextCell $ws1.Cells.Item(12, 3) 'verdana.ttf'
$ws1.Cells.Item(12, 4).Value = 14
$ws1.Cells.Item(12, 5).Value = 4
Set-TextCell $ws1.Cells.Item(12, 6) '?'

$rng1 = $ws1.Range("B12:F12")
$rng1.Copy()
$rng1.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$rng1.Style = "Normal"

# --- Sheet "Translation": update rows 68-70, add new rows 71-105 ---
$ws2 = $wb.Worksheets.Item("Translation")

# Row 68
Set-TextCell $ws2.Cells.Item(68, 2) 'SingleUseId124'
Set-TextCell $ws2.Cells.Item(68, 3) 'Small'
Set-TextCell $ws2.Cells.Item(68, 4) 'Left'
Set-TextCell $ws2.Cells.Item(68, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(68, 6) '0000'

# Row 69
Set-TextCell $ws2.Cells.Item(69, 2) 'SingleUseId125'
Set-TextCell $ws2.Cells.Item(69, 3) 'Small'
Set-TextCell $ws2.Cells.Item(69, 4) 'Left'
Set-TextCell $ws2.Cells.Item(69, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(69, 6) '0000'

# Row 70
Set-TextCell $ws2.Cells.Item(70, 2) 'SingleUseId126'
Set-TextCell $ws2.Cells.Item(70, 3) 'Small'
Set-TextCell $ws2.Cells.Item(70, 4) 'Center'
Set-TextCell $ws2.Cells.Item(70, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(70, 6) '<value>'

# Row 71
Set-TextCell $ws2.Cells.Item(71, 2) 'SingleUseId127'
Set-TextCell $ws2.Cells.Item(71, 3) 'Small'
Set-TextCell $ws2.Cells.Item(71, 4) 'Left'
Set-TextCell $ws2.Cells.Item(71, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(71, 6) '0000'

# Row 72
Set-TextCell $ws2.Cells.Item(72, 2) 'SingleUseId128'
Set-TextCell $ws2.Cells.Item(72, 3) 'Small'
Set-TextCell $ws2.Cells.Item(72, 4) 'Center'
Set-TextCell $ws2.Cells.Item(72, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(72, 6) '<value>'

# Row 73
Set-TextCell $ws2.Cells.Item(73, 2) 'SingleUseId129'
Set-TextCell $ws2.Cells.Item(73, 3) 'Small'
Set-TextCell $ws2.Cells.Item(73, 4) 'Left'
Set-TextCell $ws2.Cells.Item(73, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(73, 6) '0000'

# Row 74
Set-TextCell $ws2.Cells.Item(74, 2) 'SingleUseId130'
Set-TextCell $ws2.Cells.Item(74, 3) 'Small'
Set-TextCell $ws2.Cells.Item(74, 4) 'Center'
Set-TextCell $ws2.Cells.Item(74, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(74, 6) '<value>'

# Row 75
Set-TextCell $ws2.Cells.Item(75, 2) 'SingleUseId131'
Set-TextCell $ws2.Cells.Item(75, 3) 'Small'
Set-TextCell $ws2.Cells.Item(75, 4) 'Left'
Set-TextCell $ws2.Cells.Item(75, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(75, 6) '0000'

# Row 76
Set-TextCell $ws2.Cells.Item(76, 2) 'SingleUseId132'
Set-TextCell $ws2.Cells.Item(76, 3) 'Small'
Set-TextCell $ws2.Cells.Item(76, 4) 'Center'
Set-TextCell $ws2.Cells.Item(76, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(76, 6) '<value>'

# Row 77
Set-TextCell $ws2.Cells.Item(77, 2) 'SingleUseId133'
Set-TextCell $ws2.Cells.Item(77, 3) 'Small'
Set-TextCell $ws2.Cells.Item(77, 4) 'Left'
Set-TextCell $ws2.Cells.Item(77, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(77, 6) '0000'

# Row 78
Set-TextCell $ws2.Cells.Item(78, 2) 'SingleUseId134'
Set-TextCell $ws2.Cells.Item(78, 3) 'Small'
Set-TextCell $ws2.Cells.Item(78, 4) 'Center'
Set-TextCell $ws2.Cells.Item(78, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(78, 6) '<value>'

# Row 79
Set-TextCell $ws2.Cells.Item(79, 2) 'SingleUseId135'
Set-TextCell $ws2.Cells.Item(79, 3) 'Small'
Set-TextCell $ws2.Cells.Item(79, 4) 'Left'
Set-TextCell $ws2.Cells.Item(79, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(79, 6) '0000'

# Row 80
Set-TextCell $ws2.Cells.Item(80, 2) 'SingleUseId136'
Set-TextCell $ws2.Cells.Item(80, 3) 'Small'
Set-TextCell $ws2.Cells.Item(80, 4) 'Center'
Set-TextCell $ws2.Cells.Item(80, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(80, 6) '<value>'

# Row 81
Set-TextCell $ws2.Cells.Item(81, 2) 'SingleUseId137'
Set-TextCell $ws2.Cells.Item(81, 3) 'Small'
Set-TextCell $ws2.Cells.Item(81, 4) 'Left'
Set-TextCell $ws2.Cells.Item(81, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(81, 6) '0000'

# Row 82
Set-TextCell $ws2.Cells.Item(82, 2) 'SingleUseId138'
Set-TextCell $ws2.Cells.Item(82, 3) 'Small'
Set-TextCell $ws2.Cells.Item(82, 4) 'Center'
Set-TextCell $ws2.Cells.Item(82, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(82, 6) '<value>'

# Row 83
Set-TextCell $ws2.Cells.Item(83, 2) 'SingleUseId139'
Set-TextCell $ws2.Cells.Item(83, 3) 'Small'
Set-TextCell $ws2.Cells.Item(83, 4) 'Left'
Set-TextCell $ws2.Cells.Item(83, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(83, 6) '0000'

# Row 84
Set-TextCell $ws2.Cells.Item(84, 2) 'SingleUseId140'
Set-TextCell $ws2.Cells.Item(84, 3) 'Small'
Set-TextCell $ws2.Cells.Item(84, 4) 'Center'
Set-TextCell $ws2.Cells.Item(84, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(84, 6) '<value>'

# Row 85
Set-TextCell $ws2.Cells.Item(85, 2) 'SingleUseId141'
Set-TextCell $ws2.Cells.Item(85, 3) 'Small'
Set-TextCell $ws2.Cells.Item(85, 4) 'Left'
Set-TextCell $ws2.Cells.Item(85, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(85, 6) '0000'

# Row 86
Set-TextCell $ws2.Cells.Item(86, 2) 'SingleUseId142'
Set-TextCell $ws2.Cells.Item(86, 3) 'Small'
Set-TextCell $ws2.Cells.Item(86, 4) 'Center'
Set-TextCell $ws2.Cells.Item(86, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(86, 6) '<value>'

# Row 87
Set-TextCell $ws2.Cells.Item(87, 2) 'SingleUseId143'
Set-TextCell $ws2.Cells.Item(87, 3) 'Small'
Set-TextCell $ws2.Cells.Item(87, 4) 'Left'
Set-TextCell $ws2.Cells.Item(87, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(87, 6) '0000'

# Row 88
Set-TextCell $ws2.Cells.Item(88, 2) 'SingleUseId144'
Set-TextCell $ws2.Cells.Item(88, 3) 'Small'
Set-TextCell $ws2.Cells.Item(88, 4) 'Center'
Set-TextCell $ws2.Cells.Item(88, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(88, 6) '<value>'

# Row 89
Set-TextCell $ws2.Cells.Item(89, 2) 'SingleUseId145'
Set-TextCell $ws2.Cells.Item(89, 3) 'Small'
Set-TextCell $ws2.Cells.Item(89, 4) 'Left'
Set-TextCell $ws2.Cells.Item(89, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(89, 6) '0000'

# Row 90
Set-TextCell $ws2.Cells.Item(90, 2) 'SingleUseId146'
Set-TextCell $ws2.Cells.Item(90, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(90, 4) 'Left'
Set-TextCell $ws2.Cells.Item(90, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(90, 6) 'Parameter1'

# Row 91
Set-TextCell $ws2.Cells.Item(91, 2) 'SingleUseId147'
Set-TextCell $ws2.Cells.Item(91, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(91, 4) 'Left'
Set-TextCell $ws2.Cells.Item(91, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(91, 6) 'Parameter6'

# Row 92
Set-TextCell $ws2.Cells.Item(92, 2) 'SingleUseId148'
Set-TextCell $ws2.Cells.Item(92, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(92, 4) 'Left'
Set-TextCell $ws2.Cells.Item(92, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(92, 6) 'Parameter2'

# Row 93
Set-TextCell $ws2.Cells.Item(93, 2) 'SingleUseId149'
Set-TextCell $ws2.Cells.Item(93, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(93, 4) 'Left'
Set-TextCell $ws2.Cells.Item(93, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(93, 6) 'Parameter3'

# Row 94
Set-TextCell $ws2.Cells.Item(94, 2) 'SingleUseId150'
Set-TextCell $ws2.Cells.Item(94, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(94, 4) 'Left'
Set-TextCell $ws2.Cells.Item(94, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(94, 6) 'Parameter8'

# Row 95
Set-TextCell $ws2.Cells.Item(95, 2) 'SingleUseId151'
Set-TextCell $ws2.Cells.Item(95, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(95, 4) 'Left'
Set-TextCell $ws2.Cells.Item(95, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(95, 6) 'Parameter4'

# Row 96
Set-TextCell $ws2.Cells.Item(96, 2) 'SingleUseId152'
Set-TextCell $ws2.Cells.Item(96, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(96, 4) 'Left'
Set-TextCell $ws2.Cells.Item(96, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(96, 6) 'Parameter9'

# Row 97
Set-TextCell $ws2.Cells.Item(97, 2) 'SingleUseId153'
Set-TextCell $ws2.Cells.Item(97, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(97, 4) 'Left'
Set-TextCell $ws2.Cells.Item(97, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(97, 6) 'Parameter5'

# Row 98
Set-TextCell $ws2.Cells.Item(98, 2) 'SingleUseId154'
Set-TextCell $ws2.Cells.Item(98, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(98, 4) 'Left'
Set-TextCell $ws2.Cells.Item(98, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(98, 6) 'Parameter10'

# Row 99
Set-TextCell $ws2.Cells.Item(99, 2) 'SingleUseId155'
Set-TextCell $ws2.Cells.Item(99, 3) 'Typography_05'
Set-TextCell $ws2.Cells.Item(99, 4) 'Left'
Set-TextCell $ws2.Cells.Item(99, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(99, 6) 'Parameter7'

# Row 100
Set-TextCell $ws2.Cells.Item(100, 2) 'SingleUseId156'
Set-TextCell $ws2.Cells.Item(100, 3) 'Typography_02'
Set-TextCell $ws2.Cells.Item(100, 4) 'Left'
Set-TextCell $ws2.Cells.Item(100, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(100, 6) 'Calibration Page'

# Row 101
Set-TextCell $ws2.Cells.Item(101, 2) 'SingleUseId157'
Set-TextCell $ws2.Cells.Item(101, 3) 'Typography_01'
Set-TextCell $ws2.Cells.Item(101, 4) 'Center'
Set-TextCell $ws2.Cells.Item(101, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(101, 6) '<>'

# Row 102
Set-TextCell $ws2.Cells.Item(102, 2) 'SingleUseId158'
Set-TextCell $ws2.Cells.Item(102, 3) 'Typography_01'
Set-TextCell $ws2.Cells.Item(102, 4) 'Right'
Set-TextCell $ws2.Cells.Item(102, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(102, 6) '<>'

# Row 103
Set-TextCell $ws2.Cells.Item(103, 2) 'SingleUseId159'
Set-TextCell $ws2.Cells.Item(103, 3) 'Typography_01'
Set-TextCell $ws2.Cells.Item(103, 4) 'Center'
Set-TextCell $ws2.Cells.Item(103, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(103, 6) '<>'

# Row 104
Set-TextCell $ws2.Cells.Item(104, 2) 'SingleUseId160'
Set-TextCell $ws2.Cells.Item(104, 3) 'Typography_01'
Set-TextCell $ws2.Cells.Item(104, 4) 'Right'
Set-TextCell $ws2.Cells.Item(104, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(104, 6) '<>'

# Row 105
Set-TextCell $ws2.Cells.Item(105, 2) 'SingleUseId161'
Set-TextCell $ws2.Cells.Item(105, 3) 'Small'
Set-TextCell $ws2.Cells.Item(105, 4) 'Left'
Set-TextCell $ws2.Cells.Item(105, 5) 'LTR'
Set-TextCell $ws2.Cells.Item(105, 6) '0000'

$rng2 = $ws2.Range("B68:F105")
$rng2.Copy()
$rng2.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$rng2.Style = "Normal"